$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.348477841146678
$ws.Range("E2").Value = 22.83798409381314
$ws.Range("F2").Value = 96.96021120607674
